$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.769.83"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.088.01"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "234.10"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.08%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.626"
$c.ClearFormats()
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "58.32"
$c.ClearFormats()
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  +2.97%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "15.21"
$c.ClearFormats()
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "2.395.61"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "21.21"
$c.ClearFormats()
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("E15").Value = "  +0.98%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.36"
$c.ClearFormats()
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "2.087.47"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "37.768.54"
$ws.Range("E18").Value = "  +0.43%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.12"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.93%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "71.05"
$c.ClearFormats()
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +0.66%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "229.81"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E25").Value = "  +0.67%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.74"
$c.ClearFormats()
$ws.Range("E26").Value = "  +8.69%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "171.40"
$c.ClearFormats()
$ws.Range("E27").Value = "  +1.42%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.135"
$c.ClearFormats()
$ws.Range("E28").Value = "  -3.14%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "19.54"
$c.ClearFormats()
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +0.81%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.70"
$c.ClearFormats()
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("E38").Value = "  -0.18%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.38"
$c.ClearFormats()
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  +10.07%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "101.14"
$c.ClearFormats()
$ws.Range("E41").Value = "  +3.35%  "
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E43").Value = "  +4.17%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "16.70"
$c.ClearFormats()
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").Value = "1.452.18"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.06"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "4.10"
$c.ClearFormats()
$ws.Range("E48").Value = "  -3.52%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.ClearFormats()
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("D51").Value = "2.279.24"
$ws.Range("E51").Value = "  +0.36%  "
